# Update the "Elapsed Duration(Hrs)" (column G) values on several sheets to
# reflect the recalculated outage durations, and remove the now-empty
# trailing row 7 from sheet "R1".

$wb = $excel.ActiveWorkbook

# --- Sheet R1 ---
$ws1 = $wb.Worksheets.Item("R1")
$ws1.Range("G2").Value = "3952:27:48"
$ws1.Range("G3").Value = "92:00:26"
$ws1.Range("G4").Value = "115:00:26"
$ws1.Rows(7).Delete()

# --- Sheet R2 ---
$ws2 = $wb.Worksheets.Item("R2")
$ws2.Range("G2").Value = "12133:51:26"
$ws2.Range("G3").Value = "3263:34:55"
$ws2.Range("G4").Value = "501:46:29"

# --- Sheet R4 ---
$ws4 = $wb.Worksheets.Item("R4")
$ws4.Range("G2").Value = "2979:41:15"
$ws4.Range("G3").Value = "206:53:30"
$ws4.Range("G4").Value = "95:05:55"
$ws4.Range("G5").Value = "92:43:28"

# --- Sheet R5 ---
$ws5 = $wb.Worksheets.Item("R5")
$ws5.Range("G2").Value = "453:40:14"

# --- Sheet R6 ---
$ws6 = $wb.Worksheets.Item("R6")
$ws6.Range("G2").Value = "94:12:32"
